$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 467
$ws1.Range("F7").Value = 2544
$ws1.Range("F8").Value = 434
$ws1.Range("F9").Value = 6858
$ws1.Range("F12").Value = 4

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 467
$ws4.Range("F9").Value = 2544
$ws4.Range("F10").Value = 434
$ws4.Range("F11").Value = 6858
$ws4.Range("F14").Value = 4
$ws4.Range("F16").Value = 1
